# Add a new row (row 6) to the Income sheet, mirroring the existing rows:
# Source = "Salary", Amount = 5600, Dte = 45689.00037037037 (formatted as a date,
# same style as the rows above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text + number values for the new row
$ws.Range("A6").Value = "Salary"
$ws.Range("B6").Value = 5600

# Copy the date cell's formatting (style) from the row above so the new
# C6 cell reuses the same date number-format/style as C2:C5, then set its value.
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C6").Value = 45689.00037037037

$excel.CutCopyMode = 0
